# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" query timestamps on the existing "data" sheet
#    (panel re-queried later the same day).
# 2) Add a new "metadata" worksheet (placed after "data") describing the
#    panel query itself (name/id/version/request time/etc).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update F2:F26 ("time_taken") on the data sheet -------------------
$newTimes = @(
    "2021-10-05 14:33:45.040125",
    "2021-10-05 14:33:45.040133",
    "2021-10-05 14:33:45.040136",
    "2021-10-05 14:33:45.040139",
    "2021-10-05 14:33:45.040142",
    "2021-10-05 14:33:45.040144",
    "2021-10-05 14:33:45.040147",
    "2021-10-05 14:33:45.040149",
    "2021-10-05 14:33:45.040152",
    "2021-10-05 14:33:45.040154",
    "2021-10-05 14:33:45.040157",
    "2021-10-05 14:33:45.040159",
    "2021-10-05 14:33:45.040162",
    "2021-10-05 14:33:45.040164",
    "2021-10-05 14:33:45.040167",
    "2021-10-05 14:33:45.040169",
    "2021-10-05 14:33:45.040172",
    "2021-10-05 14:33:45.040175",
    "2021-10-05 14:33:45.040178",
    "2021-10-05 14:33:45.040180",
    "2021-10-05 14:33:45.040183",
    "2021-10-05 14:33:45.040185",
    "2021-10-05 14:33:45.040188",
    "2021-10-05 14:33:45.040190",
    "2021-10-05 14:33:45.040193"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the "metadata" worksheet, right after "data" -----------------
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Dyslipidaemia"
$ws.Range("C2").Value = 332
# data_version ("0.22") is a text label, not a number -- force text storage
# so it isn't silently coerced to the numeric value 0.22.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.22"
$ws.Range("E2").Value = "2021-06-02T23:59:37.205910Z"
$ws.Range("F2").Value = "2021-10-05 14:33:45.036363"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/332/?format=json"

$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1

$wb.Worksheets.Item("data").Activate()
